$d = $word.ActiveDocument

# 1. "Data is collected annually via snorkel surveys; video camera systems"
#    -> "Data are collected annually via snorkel surveys and a video weir. Video camera systems"
$d.Content.Find.Execute(
    "Data is collected annually via snorkel surveys; video camera systems",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Data are collected annually via snorkel surveys and a video weir. Video camera systems",
    2)

# 2. "from August-December. Data from this monitoring is used"
#    -> "from December " + Chr(8211) + " August. Data from this monitoring are used"
$d.Content.Find.Execute(
    "from August-December. Data from this monitoring is used",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ("from December " + [char]8211 + " August. Data from this monitoring are used"),
    2)

# 3. "for spring-run Chinook salmon in the Sacramento River Watershed"
#    -> "for spring-run Chinook salmon (Oncorhynchus tshawytscha) in the Sacramento River Watershed"
$d.Content.Find.Execute(
    "for spring-run Chinook salmon in the Sacramento River Watershed",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "for spring-run Chinook salmon (Oncorhynchus tshawytscha) in the Sacramento River Watershed",
    2)
